$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..30 (row 30 is the last data row, a
# full template with the correct number/style formatting for each column).
# We duplicate that row twice (inserting full rows, which correctly carries
# over cell styles such as the left-aligned "is_active" boolean column and
# the fill-flagged "email" column) to create rows 31 and 32, then overwrite
# the per-row values.

$ws.Rows.Item(30).Copy()
$ws.Rows.Item(31).Insert(-4121)   # xlShiftDown

$ws.Rows.Item(30).Copy()
$ws.Rows.Item(32).Insert(-4121)   # xlShiftDown

# Fill in row 32 (John Doe) first, then row 31 (Jane Smith), so that the
# shared-string table ends up with the same append order as the source
# workbook: John Doe, john.doe@xyz.com, Jane Smith, jane.smith@xyz.com.
$ws.Range("A32").Value = 110031
$ws.Range("B32").Value = 9317596767
$ws.Range("C32").Value = "John Doe"
$ws.Range("D32").Value = "john.doe@xyz.com"
$ws.Range("E32").Value = 818876431

$ws.Range("A31").Value = 110030
$ws.Range("B31").Value = 9317596768
$ws.Range("C31").Value = "Jane Smith"
$ws.Range("D31").Value = "jane.smith@xyz.com"
$ws.Range("E31").Value = 818876432

# Columns F, G, H, I, J, K (status_code, lang_code, last_login_method,
# is_active, cr_by, cr_dtimes) already carry the correct values and styles
# from the duplicated template row, so nothing else needs to change there.

$ws.Range("E28").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
